# Rename the sheet and test case IDs from SCD0338-014 -> SCD0026-014,
# and move the active selection/view from Q4 (with topLeftCell=J1) to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet tab name: SCD0338 -> SCD0026
$ws.Name = "SCD0026"

# TC_ID cells (column B) on the two data rows: SCD0338-014 -> SCD0026-014
$ws.Range("B2").Value = "SCD0026-014"
$ws.Range("B3").Value = "SCD0026-014"

# Reset the view/selection to B4 (also clears the old topLeftCell=J1 scroll position)
$ws.Range("B4").Select() | Out-Null

Write-Output "done"
